# Weekly update: insert a new price record as a new row 89 in the
# "Hortaliza, Feria Lagunitas de Puerto Montt - Apio" sheet, pushing the
# existing rows 89-131 down to 90-132 (dimension grows from R131 to R132).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 89; everything at/after row 89 shifts
# down by one (row 131 -> row 132).
$ws.Rows.Item(89).Insert()

# Populate the newly inserted row 89 with the new weekly record.
$ws.Cells.Item(89, 1).Value = 4
$ws.Cells.Item(89, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(89, 3).Value = "Los Lagos"
$ws.Cells.Item(89, 4).Value = 44466
$ws.Cells.Item(89, 5).Value = 10
$ws.Cells.Item(89, 6).Value = 100112017
$ws.Cells.Item(89, 7).Value = "Apio"
$ws.Cells.Item(89, 8).Value = "Americana (o)"
$ws.Cells.Item(89, 9).Value = "Primera"
$ws.Cells.Item(89, 10).Value = 30
$ws.Cells.Item(89, 11).Value = 12000
$ws.Cells.Item(89, 12).Value = 12000
$ws.Cells.Item(89, 13).Value = 12000
$ws.Cells.Item(89, 14).Value = "`$/docena de matas"
$ws.Cells.Item(89, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(89, 16).Value = 2000
$ws.Cells.Item(89, 17).Value = 6
$ws.Cells.Item(89, 18).Value = "Hortaliza"

# Make sure the date cell keeps the sheet's date style (same as D90:D132).
$ws.Cells.Item(89, 4).NumberFormat = $ws.Cells.Item(90, 4).NumberFormat

Write-Host "Inserted new row 89 (date 44466) and shifted old rows 89-131 to 90-132."
